$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Row 4: add ExploreChildSearch in D4, shift existing ChecklistName into E4
$ws.Range("D4").Value = "ExploreChildSearch"
$ws.Range("E4").Value = "ChecklistName"

# Row 5: change B5 to AUTOMATION SEARCH, D5 to Work Orders, shift old D5 value into E5
$ws.Range("B5").Value = "AUTOMATION SEARCH"
$ws.Range("D5").Value = "Work Orders"
$ws.Range("E5").Value = "SFM_Auto_2372018162541"

# Row 7: change B7 to AUTOMATION SEARCH
$ws.Range("B7").Value = "AUTOMATION SEARCH"

# Column B width change (stored OOXML width = ColumnWidth + 5/6)
$ws.Columns.Item(2).ColumnWidth = 24.5 - (5/6)

# Update selection to A14
$ws.Range("A14").Select()
